# feat: add 2022-Q1 data
#
# - Add a new quarter sheet "2022-Q1" (same look as "2021-Q4") positioned
#   right before the "总计" sheet, with that quarter's per-fund holdings.
# - Update "总计" (totals): add a new top row for "2022-Q1", keeping the
#   existing "2021-Q4"/"2021-Q3" rows below it (renumbering the leading
#   index column).
#
# Implementation note: rather than inserting a brand-new blank sheet for
# "2022-Q1", the *existing* "总计" worksheet is repurposed (renamed + its
# data replaced) to become "2022-Q1", and a brand-new sheet is appended
# right after it and renamed to "总计" with the refreshed totals. This
# keeps the workbook's internal sheet order (2021-Q3, 2021-Q4, 2022-Q1,
# 总计) while giving the *new* "总计" the newly-created sheet identity —
# matching how the workbook actually evolved.

$wb = $excel.ActiveWorkbook

$totalOld = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# New sheet for the refreshed "总计" totals, inserted right after the
# current "总计" (which is about to become "2022-Q1").
$newTotal = $wb.Worksheets.Add($null, $totalOld)

# Match the page setup / outline defaults the other worksheets in this
# workbook already use.
$newTotal.Outline.SummaryRow = 1
$newTotal.Outline.SummaryColumn = 1
$newTotal.PageSetup.FitToPagesWide = 1
$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 1) Pull formatting forward before any content gets overwritten.
# ---------------------------------------------------------------------------
# "总计"'s current header (B1:D1) and index-column (A) styling -> the new
# totals sheet.
$totalOld.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$totalOld.Range("A2").Copy()
$newTotal.Range("A2:A4").PasteSpecial(-4122)
$newTotal.Application.CutCopyMode = $false

# "2021-Q4"'s header (B1:H1) and index-column (A) styling -> the sheet
# that's about to become "2022-Q1" (needs the same 8-column fund layout).
$q4.Range("B1:H1").Copy()
$totalOld.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$totalOld.Range("A2:A3").PasteSpecial(-4122)
$totalOld.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Turn the old "总计" sheet into "2022-Q1" and fill in its fund data.
# ---------------------------------------------------------------------------
$q1 = $totalOld
$q1.Name = "2022-Q1"

# Clear any leftover totals-sheet values beyond the rows we're about to
# (re)write (old sheet had rows 1-3; new layout also uses rows 1-3, so
# nothing extra to blank out, but make sure column D's leftover label
# is replaced below).
$q1.Range("D1").Value = "基金规模"
$q1.Range("C1").Value = "基金名称"
$q1.Range("B1").Value = "基金代码"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B/D/E/F/G hold fund codes/percentages stored as plain text (not
# numbers), e.g. "006923" must keep its leading zero. A leading apostrophe
# forces text entry exactly like typing it into Excel, so it doesn't get
# auto-coerced to a number.

# Row 2: fund 006923
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'006923"
$q1.Range("C2").Value = "前海开源沪港深非周期性行业股票A"
$q1.Range("D2").Value = "'0.54"
$q1.Range("E2").Value = "'93.77"
$q1.Range("F2").Value = "'5.65"
$q1.Range("G2").Value = "'0.0305"
$q1.Range("H2").Value = 5

# Row 3: fund 006924
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'006924"
$q1.Range("C3").Value = "前海开源沪港深非周期性行业股票C"
$q1.Range("D3").Value = "'0.22"
$q1.Range("E3").Value = "'93.77"
$q1.Range("F3").Value = "'5.65"
$q1.Range("G3").Value = "'0.0124"
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 3) Populate the new "总计" sheet with the refreshed totals.
# ---------------------------------------------------------------------------
$total = $newTotal
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.26

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.03
